$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 11, pushing the existing rows 11-12
# (weekly history) down to rows 13-14.
$ws.Range("11:12").Insert()

# New row 11: latest week, "Primera" quality.
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Vega Monumental Concepción"
$ws.Range("C11").Value = "Bíobío"
$ws.Range("D11").Value = 44776
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100104
$ws.Range("H11").Value = "Frutos de pepita"
$ws.Range("I11").Value = 100104003
$ws.Range("J11").Value = "Membrillo"
$ws.Range("K11").Value = "Champion"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 10000
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 10000
$ws.Range("Q11").Value = "$/bandeja 18 kilos granel"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 556
$ws.Range("T11").Value = 18

# New row 12: latest week, "Segunda" quality.
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Vega Monumental Concepción"
$ws.Range("C12").Value = "Bíobío"
$ws.Range("D12").Value = 44776
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100104
$ws.Range("H12").Value = "Frutos de pepita"
$ws.Range("I12").Value = 100104003
$ws.Range("J12").Value = "Membrillo"
$ws.Range("K12").Value = "Champion"
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 8000
$ws.Range("O12").Value = 8000
$ws.Range("P12").Value = 8000
$ws.Range("Q12").Value = "$/bandeja 18 kilos granel"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 444
$ws.Range("T12").Value = 18

Write-Output ("Dimension after edit: " + $ws.UsedRange.Address())
